$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-14 Tuesday" "2024-05-15 Wednesday"

Replace-Text "465×4=1860" "526×5=2630"
Replace-Text "405×5=2025" "676×6=4056"
Replace-Text "151×7=1057" "880×6=5280"
Replace-Text "541×3=1623" "311×8=2488"
Replace-Text "446×5=2230" "930×5=4650"

Replace-Text "102×4=408" "975×5=4875"
Replace-Text "519×3=1557" "295×9=2655"
Replace-Text "415×6=2490" "193×2=386"
Replace-Text "304×7=2128" "555×7=3885"
Replace-Text "403×5=2015" "593×5=2965"

Replace-Text "679×9=6111" "524×9=4716"
Replace-Text "612×6=3672" "334×6=2004"
Replace-Text "106×6=636" "263×2=526"
Replace-Text "565×5=2825" "356×4=1424"
Replace-Text "120×2=240" "208×7=1456"

Replace-Text "359×6=2154" "170×3=510"
Replace-Text "588×4=2352" "400×2=800"
Replace-Text "909×9=8181" "959×7=6713"
Replace-Text "794×6=4764" "991×6=5946"
Replace-Text "531×7=3717" "998×7=6986"

Replace-Text "830×6=4980" "114×8=912"
Replace-Text "150×8=1200" "760×5=3800"
Replace-Text "838×8=6704" "952×5=4760"
Replace-Text "866×5=4330" "801×7=5607"
Replace-Text "803×2=1606" "404×8=3232"
